# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the
# 633d959d-... handoff/handback pair on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("zh-cn", "de-de")
$newHandoffTime = @{
    "zh-cn" = "2016-03-20 06:17:21"
    "de-de" = "2016-03-20 06:17:24"
}
$newHandbackTime = @{
    "zh-cn" = "2016-03-20 06:17:41"
    "de-de" = "2016-03-20 06:17:47"
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 3 and Row 4 both reference the same handoff/handback file pair,
    # so both rows carry the same timestamps.
    $ws.Range("E3").Value = $newHandoffTime[$name]
    $ws.Range("E4").Value = $newHandoffTime[$name]

    $ws.Range("H3").Value = $newHandbackTime[$name]
    $ws.Range("H4").Value = $newHandbackTime[$name]
}
